# Weekly update: shift existing Ciruela price rows down by two and
# insert the two newest weekly records (rows 14-15), per commit
# "Fruta / hortaliza, semanal".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Static columns (A,B,C,E,F,G,H,I,J) are identical for every data row
# (2..72) both before and after the edit, so just make sure the full
# A1:T72 used range is covered with them.
$ws.Range("A2:A72").Value = 7
$ws.Range("B2:B72").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C2:C72").Value = "Ñuble"
$ws.Range("E2:E72").Value = 16
$ws.Range("F2:F72").Value = "Fruta"
$ws.Range("G2:G72").Value = 100103
$ws.Range("H2:H72").Value = "Frutos de hueso (carozo)"
$ws.Range("I2:I72").Value = 100103002
$ws.Range("J2:J72").Value = "Ciruela"

# Varying columns per row, in order D,K,L,M,N,O,P,Q,R,S,T,
# for data rows 2..72 of the final layout.
$rows = @(
  @(44230, "Fortuna", "Primera", 100, 10000, 11000, 10700, "`$/caja 18 kilos granel", "Región de O'Higgins", 594, 18),
  @(44610, "Black Amber", "Primera", 60, 10000, 11000, 10500, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 583, 18),
  @(44189, "Red Beaut", "Primera", 60, 12000, 13000, 12333, "`$/caja 15 kilos granel", "Región de O'Higgins", 822, 15),
  @(44624, "Black Amber", "Primera", 120, 9500, 10000, 9750, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 542, 18),
  @(44624, "Black Amber", "Segunda", 60, 8000, 8000, 8000, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 444, 18),
  @(44186, "Red Beaut", "Primera", 120, 13000, 14000, 13417, "`$/caja 15 kilos granel", "Región de O'Higgins", 894, 15),
  @(44627, "Black Amber", "Primera", 120, 10000, 11000, 10500, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 583, 18),
  @(44627, "Black Amber", "Segunda", 60, 9000, 9000, 9000, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 500, 18),
  @(44265, "Black Amber", "Primera", 80, 9000, 10000, 9500, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 528, 18),
  @(44265, "Black Amber", "Segunda", 60, 8000, 8000, 8000, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 444, 18),
  @(44195, "Red Beaut", "Primera", 90, 13500, 14000, 13722, "`$/caja 15 kilos granel", "Región Metropolitana", 915, 15),
  @(44942, "Black Amber", "Primera", 60, 15000, 16000, 15500, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 861, 18),
  @(44949, "Black Amber", "Primera", 100, 11000, 12000, 11500, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 639, 18),
  @(44949, "Black Amber", "Segunda", 50, 9000, 9000, 9000, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 500, 18),
  @(44924, "Black Amber", "Especial", 60, 17000, 17000, 17000, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 944, 18),
  @(44924, "Black Amber", "Primera", 120, 15000, 16000, 15500, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 861, 18),
  @(44237, "Lemon", "Primera", 30, 12000, 13000, 12500, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 694, 18),
  @(44257, "Black Amber", "Primera", 68, 9000, 10000, 9559, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 531, 18),
  @(44257, "Black Amber", "Segunda", 30, 8000, 8000, 8000, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 444, 18),
  @(44267, "Angeleno", "Primera", 120, 9000, 10000, 9500, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 528, 18),
  @(44267, "Angeleno", "Segunda", 40, 8000, 8000, 8000, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 444, 18),
  @(44622, "Black Amber", "Primera", 240, 11000, 12000, 11500, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 639, 18),
  @(44622, "Black Amber", "Segunda", 80, 10000, 10000, 10000, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 556, 18),
  @(44617, "Black Amber", "Primera", 60, 10000, 10000, 10000, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 556, 18),
  @(44617, "Black Amber", "Segunda", 80, 8000, 9000, 8500, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 472, 18),
  @(44609, "Black Amber", "Primera", 100, 10000, 11000, 10500, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 583, 18),
  @(44277, "Black Amber", "Primera", 120, 9000, 10000, 9500, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 528, 18),
  @(44188, "Red Beaut", "Primera", 70, 12500, 13000, 12786, "`$/caja 15 kilos granel", "Región de O'Higgins", 852, 15),
  @(44231, "Larry Ann", "Primera", 80, 6000, 7000, 6375, "`$/bandeja 10 kilos granel", "Región Metropolitana", 638, 10),
  @(44203, "Black Amber", "Primera", 120, 9000, 10000, 9500, "`$/bandeja 10 kilos granel", "Región de O'Higgins", 950, 10),
  @(44203, "Black Amber", "Segunda", 80, 7000, 7500, 7250, "`$/caja 10 kilos", "Región de O'Higgins", 725, 10),
  @(44258, "Black Amber", "Primera", 120, 9000, 10000, 9500, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 528, 18),
  @(44620, "Black Amber", "Especial", 60, 12000, 12000, 12000, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 667, 18),
  @(44620, "Black Amber", "Primera", 120, 10000, 11000, 10500, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 583, 18),
  @(44246, "Angeleno", "Primera", 60, 10500, 11000, 10750, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 597, 18),
  @(44211, "Black Amber", "Primera", 60, 9500, 10000, 9792, "`$/caja 15 kilos granel", "Región de O'Higgins", 653, 15),
  @(44601, "Black Amber", "Primera", 60, 9000, 10000, 9500, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 528, 18),
  @(44614, "Black Amber", "Especial", 50, 13000, 13000, 13000, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 722, 18),
  @(44614, "Black Amber", "Primera", 120, 11000, 12000, 11500, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 639, 18),
  @(44202, "Black Amber", "Primera", 120, 9000, 10000, 9500, "`$/caja 10 kilos", "Provincia de Curicó", 950, 10),
  @(44931, "Black Amber", "Primera", 120, 15000, 16000, 15500, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 861, 18),
  @(44931, "Black Amber", "Segunda", 60, 13000, 13000, 13000, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 722, 18),
  @(44938, "Black Amber", "Primera", 50, 15000, 15000, 15000, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 833, 18),
  @(44938, "Black Amber", "Segunda", 50, 13000, 13000, 13000, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 722, 18),
  @(44244, "Lemon", "Primera", 60, 9000, 10000, 9500, "`$/caja 16 kilos granel", "Región de O'Higgins", 594, 16),
  @(44645, "Angeleno", "Primera", 120, 8000, 9000, 8500, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 472, 18),
  @(44645, "Angeleno", "Segunda", 60, 7000, 7000, 7000, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 389, 18),
  @(44586, "Black Amber", "Primera", 100, 9500, 10000, 9750, "`$/bandeja 18 kilos granel", "Región del Maule", 542, 18),
  @(44250, "Angeleno", "Primera", 120, 10000, 11000, 10500, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 583, 18),
  @(44215, "Black Amber", "Primera", 65, 12000, 13000, 12462, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 692, 18),
  @(44266, "Black Amber", "Primera", 120, 9000, 10000, 9500, "`$/caja 18 kilos granel", "Región de O'Higgins", 528, 18),
  @(44935, "Black Amber", "Primera", 50, 16000, 16000, 16000, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 889, 18),
  @(44935, "Black Amber", "Segunda", 30, 13000, 13000, 13000, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 722, 18),
  @(44607, "Black Amber", "Primera", 60, 11000, 12000, 11500, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 639, 18),
  @(44658, "Angeleno", "Primera", 120, 8000, 9000, 8500, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 472, 18),
  @(44658, "Angeleno", "Segunda", 60, 7000, 7000, 7000, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 389, 18),
  @(44944, "Black Amber", "Primera", 60, 14000, 15000, 14500, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 806, 18),
  @(44944, "Black Amber", "Segunda", 80, 11000, 12000, 11500, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 639, 18),
  @(44224, "Black Amber", "Especial", 100, 13000, 14000, 13500, "`$/caja 18 kilos granel", "Región de O'Higgins", 750, 18),
  @(44224, "Black Amber", "Primera", 80, 11000, 12000, 11500, "`$/caja 18 kilos granel", "Región de O'Higgins", 639, 18),
  @(44224, "Black Amber", "Segunda", 60, 10000, 10000, 10000, "`$/caja 18 kilos granel", "Región de O'Higgins", 556, 18),
  @(44223, "Black Amber", "Primera", 80, 11000, 12000, 11500, "`$/caja 18 kilos granel", "Región de O'Higgins", 639, 18),
  @(44223, "Black Amber", "Segunda", 60, 10000, 10000, 10000, "`$/caja 18 kilos granel", "Región de O'Higgins", 556, 18),
  @(44643, "Angeleno", "Primera", 100, 9000, 10000, 9500, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 528, 18),
  @(44643, "Angeleno", "Segunda", 50, 8000, 8000, 8000, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 444, 18),
  @(44259, "Black Amber", "Primera", 80, 8500, 9000, 8750, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 486, 18),
  @(44259, "Black Amber", "Segunda", 40, 8000, 8000, 8000, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 444, 18),
  @(44252, "Black Amber", "Primera", 55, 10000, 11000, 10545, "`$/caja 18 kilos empedrada", "Región de O'Higgins", 586, 18),
  @(44235, "Lemon", "Primera", 60, 13000, 14000, 13500, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 750, 18),
  @(44235, "Lemon", "Segunda", 60, 11000, 12000, 11500, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 639, 18),
  @(44595, "Black Amber", "Primera", 60, 9000, 9500, 9250, "`$/bandeja 18 kilos granel", "Provincia de Curicó", 514, 18)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
  $r = $i + 2
  $row = $rows[$i]
  $ws.Cells.Item($r, 4).Value = $row[0]
  $ws.Cells.Item($r, 11).Value = $row[1]
  $ws.Cells.Item($r, 12).Value = $row[2]
  $ws.Cells.Item($r, 13).Value = $row[3]
  $ws.Cells.Item($r, 14).Value = $row[4]
  $ws.Cells.Item($r, 15).Value = $row[5]
  $ws.Cells.Item($r, 16).Value = $row[6]
  $ws.Cells.Item($r, 17).Value = $row[7]
  $ws.Cells.Item($r, 18).Value = $row[8]
  $ws.Cells.Item($r, 19).Value = $row[9]
  $ws.Cells.Item($r, 20).Value = $row[10]
}

$ws.Range("A1").Select()
